$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.128.23"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.57%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'1.857.51"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.59%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.13%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'233.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.74%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.12%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'0.4681"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.27%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'42.81"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.29%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.2832"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -1.08%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'0.06461"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -1.59%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'20.90"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -2.84%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'0.07708"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -3.93%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'1.854.98"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -0.86%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'93.43"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -3.47%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'0.6802"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.25%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'5.058"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.92%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'264.64"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -1.24%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'30.101.76"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.75%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'13.37"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -4.09%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'0.000007553"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.80%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("E21").Value = "'  -0.12%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'2.100.71"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.79%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("E23").Value = "'  -0.11%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'5.151"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -1.96%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'6.091"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -1.63%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'9.290"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -0.92%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'165.23"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -2.06%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'18.47"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -1.91%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'1.883"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -3.11%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'1.367"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -0.17%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'0.09826"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.51%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'1.449"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -0.67%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'4.207"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -3.68%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'3.977"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -1.93%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'0.04655"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -0.57%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("E36").Value = "'  -1.62%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'0.6849"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -1.99%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("E38").Value = "'  +0.33%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("E39").Value = "'  -2.96%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'2.713"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +3.32%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'6.290"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +0.26%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'70.38"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -2.11%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'0.9998"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.12%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'0.8307"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -1.06%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'1.883"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -3.49%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'102.12"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -0.60%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'0.4046"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -2.60%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'9.111"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -0.66%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'923.31"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +1.89%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'6.920"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -1.67%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'34.06"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.89%  "
$ws.Range("E51").Style = "Normal"
